# PlanYourDegree cost-estimate workbook: add a "Takehome Amout for developer"
# column (E) that multiplies each developer's estimated hours (column B) by
# the billing-rate multiplier in $B$3, mirroring the existing "Total
# estimated Billing" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label in E5 (becomes a new shared-string entry).
$ws.Range("E5").Value = "Takehome Amout for developer"

# E6 holds the first (non-shared) formula; E7:E10 are filled with the same
# relative formula so Excel stores them as one shared-formula group, just
# like the existing D7:D10 block.
$ws.Range("E6").Formula = '=B6*$B$3'
$ws.Range("E7:E10").Formula = '=B7*$B$3'

# Widen column E to fit the new, longer header text, and widen column B to
# match column C's width now that both hold similarly long labels/values.
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668

# Move the active selection to the new column, matching the edited workbook.
$ws.Range("E11").Select()
